$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.157.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +12.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.698.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +14.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.23%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.705.95'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +14.76%  '

$ws.Range("E10").Value = '  +11.31%  '

$ws.Range("E11").Value = '  +12.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.353'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.48%  '

$ws.Range("E13").Value = '  +1.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.158.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +14.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.845.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +12.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +13.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000144'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +11.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.692.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +14.05%  '

$ws.Range("E19").Value = '  +6.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '361.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +15.03%  '

$ws.Range("E21").Value = '  +12.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.431'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.171'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.779.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0874'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +17.77%  '

$ws.Range("E30").Value = '  +8.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.84'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.17%  '

$ws.Range("E37").Value = '  +12.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.879'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.89%  '

$ws.Range("E39").Value = '  +15.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '308.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +23.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.850'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +36.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.651'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.57%  '

$ws.Range("E45").Value = '  +14.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.102'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +23.23%  '

$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +14.45%  '

$ws.Range("E50").Value = '  +9.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.057.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.46%  '
